$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.2035398230088496
$ws.Range("C2").Value = 0.5796460176991151
$ws.Range("P2").Value = 0.1592920353982301
$ws.Range("S2").Value = 0.05752212389380531

# Row 3
$ws.Range("B3").Value = 0.006993006993006993
$ws.Range("C3").Value = 0.06293706293706294
$ws.Range("J3").Value = 0.01398601398601399
$ws.Range("P3").Value = 0.7342657342657343
$ws.Range("S3").Value = 0.1818181818181818

# Row 4
$ws.Range("J4").Value = 0.02777777777777778
$ws.Range("P4").Value = 0.6388888888888888
$ws.Range("S4").Value = 0.3333333333333333

# Row 6
$ws.Range("B6").Value = 0.07222222222222222
$ws.Range("F6").Value = 0.06666666666666667
$ws.Range("J6").Value = 0.2166666666666667
$ws.Range("O6").Value = 0.01666666666666667
$ws.Range("Q6").Value = 0.1388888888888889
$ws.Range("R6").Value = 0.1111111111111111
$ws.Range("S6").Value = 0.3777777777777778

# Row 7
$ws.Range("B7").Value = 0.1224489795918367
$ws.Range("D7").Value = 0.01360544217687075
$ws.Range("E7").Value = 0.006802721088435374
$ws.Range("F7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.1020408163265306
$ws.Range("O7").Value = 0.006802721088435374
$ws.Range("Q7").Value = 0.1360544217687075
$ws.Range("R7").Value = 0.06802721088435375
$ws.Range("S7").Value = 0.4965986394557823

# Row 8
$ws.Range("B8").Value = 0.07734806629834254
$ws.Range("D8").Value = 0.0138121546961326
$ws.Range("F8").Value = 0.04696132596685083
$ws.Range("J8").Value = 0.1049723756906077
$ws.Range("O8").Value = 0.03591160220994475
$ws.Range("Q8").Value = 0.2016574585635359
$ws.Range("R8").Value = 0.08287292817679558
$ws.Range("S8").Value = 0.4364640883977901

# Row 9
$ws.Range("B9").Value = 0.05076142131979695
$ws.Range("D9").Value = 0.02030456852791878
$ws.Range("F9").Value = 0.05583756345177665
$ws.Range("J9").Value = 0.1319796954314721
$ws.Range("O9").Value = 0.02538071065989848
$ws.Range("Q9").Value = 0.1624365482233502
$ws.Range("R9").Value = 0.1015228426395939
$ws.Range("S9").Value = 0.4517766497461929

# Row 10
$ws.Range("B10").Value = 0.09677419354838709
$ws.Range("D10").Value = 0.02266782911944202
$ws.Range("F10").Value = 0.06451612903225806
$ws.Range("J10").Value = 0.1185701830863121
$ws.Range("O10").Value = 0.01394943330427201
$ws.Range("Q10").Value = 0.2388840453356582
$ws.Range("R10").Value = 0.1002615518744551
$ws.Range("S10").Value = 0.3443766346992154

# Row 11
$ws.Range("G11").Value = 0.12
$ws.Range("J11").Value = 0.08888888888888889
$ws.Range("K11").Value = 0.1822222222222222
$ws.Range("L11").Value = 0.5822222222222222
$ws.Range("S11").Value = 0.02666666666666667

# Row 12
$ws.Range("G12").Value = 0.6956521739130435
$ws.Range("J12").Value = 0.2246376811594203
$ws.Range("K12").Value = 0.01449275362318841
$ws.Range("L12").Value = 0.03623188405797102
$ws.Range("S12").Value = 0.02898550724637681

# Row 13
$ws.Range("F13").Value = 0.02325581395348837
$ws.Range("G13").Value = 0.6511627906976745
$ws.Range("J13").Value = 0.2325581395348837
$ws.Range("S13").Value = 0.09302325581395349

# Row 15
$ws.Range("F15").Value = 0.0187793427230047
$ws.Range("H15").Value = 0.1314553990610329
$ws.Range("I15").Value = 0.07511737089201878
$ws.Range("J15").Value = 0.4366197183098591
$ws.Range("K15").Value = 0.05164319248826291
$ws.Range("M15").Value = 0.004694835680751174
$ws.Range("O15").Value = 0.07511737089201878
$ws.Range("S15").Value = 0.2065727699530517

# Row 16
$ws.Range("F16").Value = 0.0125
$ws.Range("H16").Value = 0.175
$ws.Range("I16").Value = 0.05625
$ws.Range("J16").Value = 0.45
$ws.Range("K16").Value = 0.05625
$ws.Range("M16").Value = 0.0375
$ws.Range("O16").Value = 0.04375
$ws.Range("S16").Value = 0.16875

# Row 17
$ws.Range("F17").Value = 0.02158273381294964
$ws.Range("H17").Value = 0.1654676258992806
$ws.Range("I17").Value = 0.1247002398081535
$ws.Range("J17").Value = 0.4028776978417266
$ws.Range("K17").Value = 0.09592326139088729
$ws.Range("M17").Value = 0.01199040767386091
$ws.Range("O17").Value = 0.07194244604316546
$ws.Range("S17").Value = 0.105515587529976

# Row 18
$ws.Range("F18").Value = 0.005128205128205128
$ws.Range("H18").Value = 0.1641025641025641
$ws.Range("I18").Value = 0.09743589743589744
$ws.Range("J18").Value = 0.4871794871794872
$ws.Range("K18").Value = 0.07179487179487179
$ws.Range("M18").Value = 0.02564102564102564
$ws.Range("O18").Value = 0.06666666666666667
$ws.Range("S18").Value = 0.08205128205128205

# Row 19
$ws.Range("F19").Value = 0.01559633027522936
$ws.Range("H19").Value = 0.1935779816513762
$ws.Range("I19").Value = 0.09357798165137615
$ws.Range("J19").Value = 0.3798165137614679
$ws.Range("K19").Value = 0.09724770642201835
$ws.Range("M19").Value = 0.02568807339449541
$ws.Range("N19").Value = 0.0009174311926605505
$ws.Range("O19").Value = 0.07889908256880734
$ws.Range("S19").Value = 0.1146788990825688
